# The underlying XML diff for this revision is purely a re-serialization
# artifact: every changed part (document.xml, endnotes.xml, footer1.xml,
# footnotes.xml, header1.xml, styles.xml, theme1.xml) differs *only* in the
# ordering of the xmlns:* namespace-prefix declarations on each part's root
# element (e.g. xmlns:r/xmlns:w15/xmlns:w14 vs. xmlns:m/xmlns:w14/xmlns:r).
# Attribute order carries no semantic meaning in XML, and no actual
# document content, text, formatting, styles, headers/footers, or theme
# data changed at all between the two revisions (the commit - "Rework DOCX
# implementation to better support testing and decrease coupling" - simply
# swapped the internal OOXML writer that re-saved the fixture).
#
# There is therefore no Word object-model action to perform: the visible
# document is identical before and after. We touch ActiveDocument so the
# session is exercised, but make no content-altering calls.
$d = $word.ActiveDocument
$null = $d.Content
